$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 9187
$ws.Range("E2").Value = 497
$ws.Range("F2").Value = 497
$ws.Range("G2").Value = 466
$ws.Range("H2").Value = 10815
$ws.Range("I2").Value = 10787
$ws.Range("J2").Value = 28
$ws.Range("K2").Value = 16794
$ws.Range("L2").Value = 7343
$ws.Range("M2").Value = 9452
$ws.Range("N2").Value = 9165
$ws.Range("O2").Value = 287
$ws.Range("P2").Value = 551
$ws.Range("Q2").Value = 2065
$ws.Range("R2").Value = -5087
$ws.Range("S2").Value = -2470
$ws.Range("T2").Value = 1625
$ws.Range("U2").Value = 439
$ws.Range("V2").Value = 5348
$ws.Range("W2").Value = 5.41
$ws.Range("X2").Value = 117.73
$ws.Range("Y2").Value = 83.96
$ws.Range("Z2").Value = 32.01
$ws.Range("AA2").Value = 77.68
$ws.Range("AB2").Value = 4789.13
$ws.Range("AC2").Value = 72030
$ws.Range("AD2").Value = 0.96
$ws.Range("AE2").Value = 85901
$ws.Range("AF2").Value = 0.81
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 0.72
$ws.Range("AI2").Value = 0.49
$ws.Range("AJ2").Value = 10802691

# Row 3
$ws.Range("D3").Value = 8921
$ws.Range("E3").Value = 960
$ws.Range("F3").Value = 960
$ws.Range("G3").Value = 924
$ws.Range("H3").Value = 797
$ws.Range("I3").Value = 752
$ws.Range("J3").Value = 45
$ws.Range("K3").Value = 16060
$ws.Range("L3").Value = 6633
$ws.Range("M3").Value = 9427
$ws.Range("N3").Value = 9097
$ws.Range("O3").Value = 330
$ws.Range("P3").Value = 551
$ws.Range("Q3").Value = 408
$ws.Range("R3").Value = -60
$ws.Range("S3").Value = -530
$ws.Range("T3").Value = 187
$ws.Range("U3").Value = 221
$ws.Range("V3").Value = 4882
$ws.Range("W3").Value = 10.76
$ws.Range("X3").Value = 8.93
$ws.Range("Y3").Value = 8.23
$ws.Range("Z3").Value = 4.85
$ws.Range("AA3").Value = 70.37
$ws.Range("AB3").Value = 1747.03
$ws.Range("AC3").Value = 6959
$ws.Range("AD3").Value = 9.25
$ws.Range("AE3").Value = 85265
$ws.Range("AF3").Value = 0.76
$ws.Range("AG3").Value = 1200
$ws.Range("AH3").Value = 1.86
$ws.Range("AI3").Value = 17.03
$ws.Range("AJ3").Value = 10802691

# Row 4
$ws.Range("D4").Value = 8445
$ws.Range("E4").Value = 935
$ws.Range("F4").Value = 1159
$ws.Range("G4").Value = 784
$ws.Range("H4").Value = 804
$ws.Range("I4").Value = 753
$ws.Range("J4").Value = 51
$ws.Range("K4").Value = 19140
$ws.Range("L4").Value = 9575
$ws.Range("M4").Value = 9565
$ws.Range("N4").Value = 9196
$ws.Range("O4").Value = 369
$ws.Range("P4").Value = 551
$ws.Range("Q4").Value = 859
$ws.Range("R4").Value = -105
$ws.Range("S4").Value = -1214
$ws.Range("T4").Value = 131
$ws.Range("U4").Value = 728
$ws.Range("V4").Value = 6850
$ws.Range("W4").Value = 11.07
$ws.Range("X4").Value = 9.52
$ws.Range("Y4").Value = 8.24
$ws.Range("Z4").Value = 4.57
$ws.Range("AA4").Value = 100.11
$ws.Range("AB4").Value = 1676.3
$ws.Range("AC4").Value = 6975
$ws.Range("AD4").Value = 8.69
$ws.Range("AE4").Value = 86195
$ws.Range("AF4").Value = 0.7
$ws.Range("AG4").Value = 1250
$ws.Range("AH4").Value = 2.06
$ws.Range("AI4").Value = 17.7
$ws.Range("AJ4").Value = 10802691

# Row 5
$ws.Range("D5").Value = 8907
$ws.Range("E5").Value = 637
$ws.Range("F5").Value = 637
$ws.Range("G5").Value = 764
$ws.Range("H5").Value = 435
$ws.Range("I5").Value = 430
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 18523
$ws.Range("L5").Value = 9186
$ws.Range("M5").Value = 9337
$ws.Range("N5").Value = 9337
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 551
$ws.Range("Q5").Value = 526
$ws.Range("R5").Value = 20
$ws.Range("S5").Value = -297
$ws.Range("T5").Value = 275
$ws.Range("U5").Value = 250
$ws.Range("V5").Value = 6436
$ws.Range("W5").Value = 7.15
$ws.Range("X5").Value = 4.89
$ws.Range("Y5").Value = 4.64
$ws.Range("Z5").Value = 2.31
$ws.Range("AA5").Value = 98.38
$ws.Range("AB5").Value = 1739.11
$ws.Range("AC5").Value = 3982
$ws.Range("AD5").Value = 15.75
$ws.Range("AE5").Value = 87515
$ws.Range("AF5").Value = 0.72
$ws.Range("AG5").Value = 1350
$ws.Range("AH5").Value = 2.15
$ws.Range("AI5").Value = 33.48
$ws.Range("AJ5").Value = 10802691

# Row 6
$ws.Range("D6").Value = 8865
$ws.Range("E6").Value = 575
$ws.Range("F6").Value = 575
$ws.Range("G6").Value = 352
$ws.Range("H6").Value = 131
$ws.Range("I6").Value = 131
$ws.Range("K6").Value = 17241
$ws.Range("L6").Value = 8589
$ws.Range("M6").Value = 8652
$ws.Range("N6").Value = 8652
$ws.Range("P6").Value = 551
$ws.Range("Q6").Value = 8
$ws.Range("R6").Value = 372
$ws.Range("S6").Value = -772
$ws.Range("T6").Value = 278
$ws.Range("U6").Value = -270
$ws.Range("V6").Value = 5896
$ws.Range("W6").Value = 6.49
$ws.Range("X6").Value = 1.48
$ws.Range("Y6").Value = 1.46
$ws.Range("Z6").Value = 0.73
$ws.Range("AA6").Value = 99.27
$ws.Range("AB6").Value = 1640.99
$ws.Range("AC6").Value = 1214
$ws.Range("AD6").Value = 36.05
$ws.Range("AE6").Value = 82464
$ws.Range("AF6").Value = 0.53
$ws.Range("AG6").Value = 2000
$ws.Range("AH6").Value = 4.57
$ws.Range("AI6").Value = 160.41
$ws.Range("AJ6").Value = 10802691

# Row 7
$ws.Range("D7").Value = 8033
$ws.Range("E7").Value = 792
$ws.Range("G7").Value = 548
$ws.Range("H7").Value = 413
$ws.Range("I7").Value = 413
$ws.Range("K7").Value = 17836
$ws.Range("L7").Value = 8939
$ws.Range("M7").Value = 8896
$ws.Range("N7").Value = 8896
$ws.Range("P7").Value = 550
$ws.Range("Q7").Value = 622
$ws.Range("R7").Value = 273
$ws.Range("S7").Value = -37
$ws.Range("T7").Value = 197
$ws.Range("U7").Value = 265
$ws.Range("W7").Value = 9.86
$ws.Range("X7").Value = 5.15
$ws.Range("Y7").Value = 4.71
$ws.Range("Z7").Value = 2.36
$ws.Range("AA7").Value = 100.48
$ws.Range("AC7").Value = 3944
$ws.Range("AD7").Value = 10.6
$ws.Range("AE7").Value = 84953
$ws.Range("AF7").Value = 0.49
$ws.Range("AG7").Value = 1875
$ws.Range("AH7").Value = 4.49
$ws.Range("AI7").Value = 47.5

# Row 8
$ws.Range("D8").Value = 8404
$ws.Range("E8").Value = 872
$ws.Range("G8").Value = 637
$ws.Range("H8").Value = 483
$ws.Range("I8").Value = 484
$ws.Range("K8").Value = 18188
$ws.Range("L8").Value = 8988
$ws.Range("M8").Value = 9200
$ws.Range("N8").Value = 9200
$ws.Range("P8").Value = 550
$ws.Range("Q8").Value = 630
$ws.Range("R8").Value = -209
$ws.Range("S8").Value = -270
$ws.Range("T8").Value = 183
$ws.Range("U8").Value = 348
$ws.Range("W8").Value = 10.38
$ws.Range("X8").Value = 5.75
$ws.Range("Y8").Value = 5.34
$ws.Range("Z8").Value = 2.68
$ws.Range("AA8").Value = 97.7
$ws.Range("AC8").Value = 4617
$ws.Range("AD8").Value = 9.05
$ws.Range("AE8").Value = 87853
$ws.Range("AF8").Value = 0.48
$ws.Range("AG8").Value = 1883
$ws.Range("AH8").Value = 4.51
$ws.Range("AI8").Value = 40.79

# Row 9
$ws.Range("D9").Value = 8723
$ws.Range("E9").Value = 929
$ws.Range("G9").Value = 704
$ws.Range("H9").Value = 532
$ws.Range("I9").Value = 534
$ws.Range("K9").Value = 18508
$ws.Range("L9").Value = 8954
$ws.Range("M9").Value = 9552
$ws.Range("N9").Value = 9554
$ws.Range("P9").Value = 550
$ws.Range("Q9").Value = 671
$ws.Range("R9").Value = -234
$ws.Range("S9").Value = -308
$ws.Range("T9").Value = 222
$ws.Range("U9").Value = 349
$ws.Range("W9").Value = 10.65
$ws.Range("X9").Value = 6.1
$ws.Range("Y9").Value = 5.7
$ws.Range("Z9").Value = 2.9
$ws.Range("AA9").Value = 93.73
$ws.Range("AC9").Value = 5101
$ws.Range("AD9").Value = 8.19
$ws.Range("AE9").Value = 91238
$ws.Range("AF9").Value = 0.46
$ws.Range("AG9").Value = 1892
$ws.Range("AH9").Value = 4.53
$ws.Range("AI9").Value = 37.09
